$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Rv0642c
$ws.Range("A3").Value = "Rv0642c"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "mmaA4 hma mma4 Rv0642c"
$ws.Range("D3").Value = "FUNCTION: Involved in the biosynthesis of hydroxymycolate, a common precursor of oxygenated mycolic acids (methoxy-mycolate and keto-mycolate). Probably transfers a methyl group from the S-adenosylmethionine (SAM) cofactor and, subsequently or simultaneously, a water molecule onto the double bound of ethylene substrates, leading to the formation of the hydroxylated product at the distal position. Involved in the activation of the antitubercular drug thiacetazone (TAC). {ECO:0000269|PubMed:10844652, ECO:0000269|PubMed:12473649}."
$ws.Range("E3").Value = 8

# Row 8: Rv2224c
$ws.Range("A8").Value = "Rv2224c"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "hip1 caeA Rv2224c MTCY427.05c"
$ws.Range("D8").Value = "FUNCTION: Serine protease that promotes tuberculosis (TB) pathogenesis by promoting the processing and the extracellular release of the M.tuberculosis (Mtb) heat-shock protein GroEL2 (PubMed:18172199, PubMed:24830429, PubMed:28346784). Hip1-dependent cleavage of multimeric GroEL2 results in release of cleaved monomeric GroEL2 into the extracellular milieu. Conversion of multimeric GroEL2 into monomeric GroEL2 is likely to be a mechanism for regulating GroEL2 functions during Mtb pathogenesis (PubMed:24830429). In vitro, exhibits proteolytic activity against synthetic peptides and the general protease substrate azocasein, and exhibits esterase activity against the ester substrate p-nitrophenylbutyrate (PubMed:24830429, PubMed:28346784). {ECO:0000269|PubMed:18172199, ECO:0000269|PubMed:24830429, ECO:0000269|PubMed:28346784}.; FUNCTION: Key immunomodulatory virulence factor, which promotes survival in host macrophages and modulates host immune responses (PubMed:18172199, PubMed:21947769, PubMed:24659689). Impacts host innate immune responses by preventing robust macrophage activation (PubMed:18172199, PubMed:21947769). Dampens macrophage proinflammatory responses by limiting toll-like receptor 2 (TLR2) activation. It also dampens TLR2-independent activation of the inflammasome and limits secretion of interleukin-18 (IL-18). May act by masking cell surface interactions between TLR2 agonists on Mtb and TLR2 on macrophages (PubMed:21947769). In addition, impacts host adaptive immune responses. It prevents robust maturation of infected dendritic cells (DCs), limits the secretion of key proinflammatory cytokines such as IL-12, impairs Ag presentation, and modulates the nature of Ag-specific T-cell responses (PubMed:24659689). {ECO:0000269|PubMed:18172199, ECO:0000269|PubMed:21947769, ECO:0000269|PubMed:24659689}."
$ws.Range("E8").Value = 8

# Row 9: Rv1410c
$ws.Range("A9").Value = "Rv1410c"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = "Rv1410c"
$ws.Range("D9").Value = "FUNCTION: In association with lipoprotein LprG probably transports triacylglycerides (TAG) across the inner cell membrane into the periplasm; TAG probably regulates lipid metabolism and growth regulation (PubMed:26751071). Confers resistance to ethidium bromide, possibly acting as an efflux pump, requires LprG lipoprotein for normal function (PubMed:18156250). With LprG maintains cell wall permeability (PubMed:21762531). Probably required with LprG for normal surface localization of LAM (PubMed:25232742, PubMed:25356793). Overexpression of LprG and Rv1410c leads to increased levels of TAG in the culture medium (PubMed:26751071). {ECO:0000269|PubMed:18156250, ECO:0000269|PubMed:21762531, ECO:0000269|PubMed:26751071, ECO:0000305|PubMed:25232742, ECO:0000305|PubMed:25356793}."
$ws.Range("E9").Value = 8

# Row 10: Rv0129c
$ws.Range("A10").Value = "Rv0129c"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "fbpC mpt45 Rv0129c MTCI5.03c"
$ws.Range("D10").Value = "FUNCTION: The antigen 85 proteins (FbpA, FbpB, FbpC) are responsible for the high affinity of mycobacteria to fibronectin, a large adhesive glycoprotein, which facilitates the attachment of M.tuberculosis to murine alveolar macrophages (AMs). They also help to maintain the integrity of the cell wall by catalyzing the transfer of mycolic acids to cell wall arabinogalactan and through the synthesis of alpha,alpha-trehalose dimycolate (TDM, cord factor). They catalyze the transfer of a mycoloyl residue from one molecule of alpha,alpha-trehalose monomycolate (TMM) to another TMM, leading to the formation of TDM. {ECO:0000269|PubMed:1830294, ECO:0000269|PubMed:9162010}."
$ws.Range("E10").Value = 8

# Row 14: Rv1821
$ws.Range("A14").Value = "Rv1821"
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "secA2 Rv1821 MTCY1A11.22c"
$ws.Range("D14").Value = "FUNCTION: Part of the Sec protein translocase complex. Interacts with the SecYEG preprotein conducting channel. Has a central role in coupling the hydrolysis of ATP to the transfer of proteins into and across the cell membrane, serving as an ATP-driven molecular motor driving the stepwise translocation of polypeptide chains across the membrane. {ECO:0000255|HAMAP-Rule:MF_01382}."
$ws.Range("E14").Value = 8

# Row 15: Rv3794
$ws.Range("A15").Value = "Rv3794"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "embA Rv3794 MTCY13D12.28"
$ws.Range("D15").Value = "FUNCTION: Arabinosyl transferase responsible for the polymerization of arabinose into the arabinan of arabinogalactan."
$ws.Range("E15").Value = 8

# Row 16: Rv1698
$ws.Range("A16").Value = "Rv1698"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = "mctB Rv1698 MTCI125.20"
$ws.Range("D16").Value = "FUNCTION: Pore-forming protein, which is involved in efflux of copper across the outer membrane. Essential for copper resistance and maintenance of a low intracellular copper concentration. Required for virulence. {ECO:0000269|PubMed:18434314, ECO:0000269|PubMed:21205886}."
$ws.Range("E16").Value = 8

# Row 17: Rv1512
$ws.Range("A17").Value = "Rv1512"
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = "epiA fcl Rv1512"
$ws.Range("D17").Value = "FUNCTION: Catalyzes the two-step NADP-dependent conversion of GDP-4-dehydro-6-deoxy-D-mannose to GDP-fucose, involving an epimerase and a reductase reaction. {ECO:0000256|HAMAP-Rule:MF_00956}."
$ws.Range("E17").Value = 8

# Row 18: Rv1433
$ws.Range("A18").Value = "Rv1433"
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = "Rv1433 RVBD_1433 P425_01489"
$ws.Range("D18").Value = "FUNCTION: Probable L,D-transpeptidase that may perform as-yet-unknown cross-linking reactions in M.tuberculosis. Is not able to generate 3->3 cross-links in peptidoglycan, using tetrapeptide stems as acyl donor substrates. May function in the anchoring of proteins to peptidoglycan. {ECO:0000269|PubMed:24041897}."
$ws.Range("E18").Value = 8

# Row 19: Rv0994
$ws.Range("A19").Value = "Rv0994"
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "moeA1 moeA Rv0994 MTCI237.08"
$ws.Range("D19").Value = "FUNCTION: Catalyzes the insertion of molybdate into adenylated molybdopterin with the concomitant release of AMP. {ECO:0000250}."
$ws.Range("E19").Value = 8

# Row 20: Rv2165c
$ws.Range("A20").Value = "Rv2165c"
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "rsmH mraW Rv2165c MTCY270.03"
$ws.Range("D20").Value = "FUNCTION: Specifically methylates the N4 position of cytidine in position 1402 (C1402) of 16S rRNA. {ECO:0000255|HAMAP-Rule:MF_01007}."
$ws.Range("E20").Value = 8

# Row 21: Rv0436c
$ws.Range("A21").Value = "Rv0436c"
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "pssA Rv0436c MTCY22G10.33c"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = 8

# Row 22: Rv0111
$ws.Range("A22").Value = "Rv0111"
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "Rv0111"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = 8

# Row 23: Rv0472c
$ws.Range("A23").Value = "Rv0472c"
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "Rv0472c"
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = 8

# Row 24: Rv3779
$ws.Range("A24").Value = "Rv3779"
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "Rv3779"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = 8

# Row 25: Rv0708
$ws.Range("A25").Value = "Rv0708"
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "rplP Rv0708 MTCY210.27"
$ws.Range("D25").Value = "FUNCTION: Binds 23S rRNA and is also seen to make contacts with the A and possibly P site tRNAs. {ECO:0000255|HAMAP-Rule:MF_01342}."
$ws.Range("E25").Value = 8

# Row 26: Rv1435c
$ws.Range("A26").Value = "Rv1435c"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Rv1435c"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = 8

# Row 27: Rv3631
$ws.Range("A27").Value = "Rv3631"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Rv3631"
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = 8

# Row 28: Rv0999
$ws.Range("A28").Value = "Rv0999"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "Rv0999"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = 8

# Row 29: Rv3632
$ws.Range("A29").Value = "Rv3632"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = "Rv3632"
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = 8

# Row 30: Rv1244
$ws.Range("A30").Value = "Rv1244"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "lpqZ Rv1244"
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = 8

# Row 31: Rv3005c
$ws.Range("A31").Value = "Rv3005c"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = "Rv3005c"
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = 8

# Row 32: Rv3267
$ws.Range("A32").Value = "Rv3267"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "Rv3267"
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = 8

# Row 33: Rv0179c
$ws.Range("A33").Value = "Rv0179c"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "lprO Rv0179c"
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = 8

# Row 34: Rv2169c
$ws.Range("A34").Value = "Rv2169c"
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = "Rv2169c"
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = 8

# Row 35: Rv0049
$ws.Range("A35").Value = "Rv0049"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "Rv0049 MTCY21D4.12"
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = 8

# Row 36: Rv0204c
$ws.Range("A36").Value = "Rv0204c"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = "Rv0204c"
$ws.Range("D36").Value = ""
$ws.Range("E36").Value = 8

